$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update each nha vuon row with its own unique code (NV_2 .. NV_9).
# B5 already holds NV_1 and is left untouched.
$ws.Range("B6").Value = "NV_2"
$ws.Range("B7").Value = "NV_3"
$ws.Range("B8").Value = "NV_4"
$ws.Range("B9").Value = "NV_5"
$ws.Range("B10").Value = "NV_6"
$ws.Range("B11").Value = "NV_7"
$ws.Range("B12").Value = "NV_8"
$ws.Range("B13").Value = "NV_9"

# Move the active selection, matching the author's last cursor position.
$ws.Range("I20").Select() | Out-Null
